$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Register the "No Spacing" built-in paragraph style (styleId
#    "NoSpacing") the same way Word does the first time the style is
#    applied to a paragraph - this materialises the <w:style> entry in
#    styles.xml.
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$seedPara = $endRange.Paragraphs.Add()
$seedPara.Range.Style = "No Spacing"
$seedPara.Range.Text = ""

$noSpacing = $d.Styles("No Spacing")
$noSpacing.Priority = 1
$pf = $noSpacing.ParagraphFormat
$pf.SpaceAfter = 0
$pf.LineSpacingRule = 0

# remove the scratch paragraph used only to seed the style definition
$seedPara.Range.Delete()

# ---------------------------------------------------------------------
# 2. Append the new notes section after the "Follow" paragraph, at the
#    end of the document body, using the exact OOXML for each
#    paragraph so styles/formatting match precisely.
# ---------------------------------------------------------------------
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$fragment = @"
<w:p $w/>
<w:p $w><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>New Notes &#8211; Wireframe discussion</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Change from beards, to albums</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Need to see full image</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>View others profiles + follow system</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>
<w:p $w><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>
"@

$insertRange = $d.Content
$insertRange.Collapse(0)
$insertRange.InsertXML($fragment)
